# Adds the new observation row (row 6) to the "Artfynd" sheet, matching the
# target diff: a Talltita (Poecile montanus) sighting reported by Anette
# Källman. The sheet's used range grows from A1:AY5 to A1:AY6 automatically
# once row 6 gets data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Id / taxon info ---------------------------------------------------
$ws.Cells.Item(6,1).Value  = 131191459          # A  Id
$ws.Cells.Item(6,2).Value  = 58043              # B  Taxonsorteringsordning
# C Valideringsstatus -> left blank (not present in source row)
$ws.Cells.Item(6,4).Value  = "NT"               # D  Rödlistade
$ws.Cells.Item(6,5).Value  = 103021             # E  TaxonId
$ws.Cells.Item(6,6).Value  = "Talltita"         # F  Artnamn
$ws.Cells.Item(6,7).Value  = "Poecile montanus" # G  Vetenskapligt namn
$ws.Cells.Item(6,8).Value  = "(Conrad von Baldenstein, 1827)" # H Auktor

# Columns that are present in the row but hold an *empty* string in the
# source data. Plain `.Value = ""` clears/omits a cell entirely in this
# engine (same as real Excel - empty input means "no content"), so instead
# we touch a boolean font property with its already-default value, which
# is enough to materialise the cell without allocating a new style.
$ws.Cells.Item(6,9).Font.Bold  = $false         # I  Antal
$ws.Cells.Item(6,11).Font.Bold = $false         # K  Ålder-Stadium
$ws.Cells.Item(6,12).Font.Bold = $false         # L  Kön

$ws.Cells.Item(6,13).Value = "lockläte, övriga läten" # M Aktivitet

$ws.Cells.Item(6,14).Font.Bold = $false         # N  Metod

# --- Location -----------------------------------------------------------
$ws.Cells.Item(6,16).Value = "Sjöberga 1:2, Ög" # P  Lokalnamn
$ws.Cells.Item(6,17).Value = 567527             # Q  Ost
$ws.Cells.Item(6,18).Value = 6509582            # R  Nord
$ws.Cells.Item(6,19).Value = 10                 # S  Noggrannhet
$ws.Cells.Item(6,20).Value = "Östergötland"      # T  Län
$ws.Cells.Item(6,21).Value = "Norrköping"       # U  Kommun
$ws.Cells.Item(6,22).Value = "Östergötland"      # V  Provins
$ws.Cells.Item(6,23).Value = "Kvillinge"        # W  Socken

# --- Dates (must stay literal text, not auto-convert to a date serial) --
# Formatting the cell as Text before assignment stops Excel's date
# autodetection; re-applying the builtin "Normal" style afterwards drops
# the temporary Text number-format override so the saved cell keeps the
# workbook's original (default) style, matching the target exactly.
$ws.Cells.Item(6,25).NumberFormat = "@"
$ws.Cells.Item(6,25).Value = "2026-02-16"       # Y  Startdatum
$ws.Cells.Item(6,25).Style = "Normal"

$ws.Cells.Item(6,27).NumberFormat = "@"
$ws.Cells.Item(6,27).Value = "2026-02-16"       # AA Slutdatum
$ws.Cells.Item(6,27).Style = "Normal"

# --- Flags ---------------------------------------------------------------
$ws.Cells.Item(6,30).Value = $false             # AD Ej återfunnen
$ws.Cells.Item(6,31).Value = $false             # AE Osäker artbestämning
$ws.Cells.Item(6,33).Value = $false             # AG Ospontan

$ws.Cells.Item(6,46).Font.Bold = $false         # AT Bestämningsår (empty)

# --- Reporter -------------------------------------------------------------
$ws.Cells.Item(6,49).Value = "Anette Källman"   # AW Rapportör
$ws.Cells.Item(6,50).Value = "Anette Källman"   # AX Observatörer

$ws.Cells.Item(6,51).Font.Bold = $false         # AY Projektnamn (empty)
